# "fix insert new column before first to work in linux"
#
# The original bug: inserting a brand new column before column A (the
# first column) and removing it again left the sheet's column metadata
# in a state where the three data columns (Date / TemperatureC /
# Summary) no longer carried an explicit, correct width - on Linux the
# auto "best fit" width that Windows Excel had cached no longer applied
# cleanly. The fix re-establishes an explicit custom width for every
# column touched by the insert/delete round trip and drops the stale
# "best fit" flag that only made sense for the old, Windows-computed
# width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reproduce the edit that exposed the bug: insert a new column before
# the first column, then remove it again so the data layout (A:C) is
# unchanged but the column-width bookkeeping gets recomputed.
$ws.Columns.Item(1).Insert() | Out-Null
$ws.Columns.Item(1).Delete() | Out-Null

# Re-apply explicit widths to the three visible columns so they no
# longer rely on the (platform-dependent) "best fit" auto width.
$ws.Range("A1").ColumnWidth = 10.833333333333332
$ws.Range("B1").ColumnWidth = 15
$ws.Range("C1").ColumnWidth = 16.833333333333336
